# Build the "Testing file" header block (merged Q12:Y23, centered)
# and the small list below it (Q26:Q33), matching the target worksheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Center-align the block first, then set the title text, then merge —
# this ordering reproduces the expected single extra (centered) cell style.
$ws.Range("Q12:Y23").HorizontalAlignment = -4108   # xlCenter
$ws.Range("Q12").Value = "Testing file"
$ws.Range("Q12:Y23").Merge()

# Small list of values underneath the merged block.
$ws.Range("Q26").Value = "a"
$ws.Range("Q27").Value = "sd"
$ws.Range("Q28").Value = "ad"
$ws.Range("Q29").Value = "ads"
$ws.Range("Q30").Value = "das"
$ws.Range("Q31").Value = "dsa"
$ws.Range("Q32").Value = "dsa"
$ws.Range("Q33").Value = "asdsda"

# Match the saved selection/active cell from the source file.
[void]$ws.Range("Q34").Select()
